# Replace the logged job-search activity rows (2-5) with a new set of
# entries, as part of testing the automation process.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Module code names picked up now that the workbook is being driven by the
# automation script (VBA project identifiers).
$wb.CodeName = "ThisWorkbook"
$ws.CodeName = "Sheet1"

# Row 2
$ws.Range("A2").Value = "<2018-11-18 Sun 21:56>"
$ws.Range("B2").Value = "Systems Integrator VAR"
$ws.Range("C2").Value = "`$90K/yr"
$ws.Range("D2").Value = "http://www.linkedin.com"
$ws.Range("E2").Value = "Matthew Handler  - Recruitment Consultant"
$ws.Range("F2").Value = "web"
$ws.Range("G2").Value = "Connected"

# Row 3
$ws.Range("A3").Value = "<2018-11-19 Mon 21:56>"
$ws.Range("B3").Value = "Systems Integrator VAR"
$ws.Range("C3").Value = "`$90K/yr"
$ws.Range("D3").Value = "http://www.linkedin.com"
$ws.Range("E3").Value = "Sankat Arbat  - Recruiter"
$ws.Range("F3").Value = "web"
$ws.Range("G3").Value = "Connected"

# Row 4
$ws.Range("A4").Value = "<2018-11-20 Tue 21:56>"
$ws.Range("B4").Value = "Systems Integrator VAR"
$ws.Range("C4").Value = "`$90K/yr"
$ws.Range("D4").Value = "http://www.linkedin.com"
$ws.Range("E4").Value = "Josh Leventhal  - Marketer/Writer"
$ws.Range("F4").Value = "web"
$ws.Range("G4").Value = "Connected"

# Row 5
$ws.Range("A5").Value = "<2018-11-21 Wed 21:56>"
$ws.Range("B5").Value = "Systems Integrator VAR"
$ws.Range("C5").Value = "`$90K/yr"
$ws.Range("D5").Value = "http://www.linkedin.com"
$ws.Range("E5").Value = "MIRZA ASLAMULLAH BAIG - DBA (Database Administrator)"
$ws.Range("F5").Value = "web"
$ws.Range("G5").Value = "Connected"
